$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Wilstermann vs Tomayapo) odds updates
$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3.4
$ws.Range("T2").Value = 8
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 9
$ws.Range("W2").Value = 21
$ws.Range("X2").Value = 17
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 17
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 29

# Row 5 (Kuressaare vs Flora) odds updates
$ws.Range("G5").Value = 7.5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 1.29
$ws.Range("T5").Value = 20
$ws.Range("V5").Value = 19.5
$ws.Range("W5").Value = 120
$ws.Range("X5").Value = 60
$ws.Range("Y5").Value = 50
$ws.Range("Z5").Value = 16.5
$ws.Range("AA5").Value = 9.25
$ws.Range("AB5").Value = 16
$ws.Range("AC5").Value = 60
$ws.Range("AD5").Value = 350
$ws.Range("AE5").Value = 7.6
$ws.Range("AF5").Value = 6.2
$ws.Range("AG5").Value = 7.5
$ws.Range("AH5").Value = 7.2
$ws.Range("AJ5").Value = 19

# Row 7 (Narva vs Harju JK Laagri) odds updates
$ws.Range("G7").Value = 1.6
$ws.Range("I7").Value = 4.6
$ws.Range("T7").Value = 8.5
$ws.Range("U7").Value = 8.5
$ws.Range("V7").Value = 6.9
$ws.Range("W7").Value = 11.5
$ws.Range("Y7").Value = 14.5
$ws.Range("Z7").Value = 15
$ws.Range("AA7").Value = 6.9
$ws.Range("AB7").Value = 10.75
$ws.Range("AC7").Value = 32
$ws.Range("AE7").Value = 14
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 65
$ws.Range("AI7").Value = 32
$ws.Range("AJ7").Value = 28

# Row 8 (Braunschweig vs Saarbrucken) odds updates
$ws.Range("L8").Value = 1.29
$ws.Range("M8").Value = 3.5
$ws.Range("N8").Value = 1.9
$ws.Range("O8").Value = 1.9

# Row 10 (Daejeon vs Pohang) odds updates
$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 3.2
$ws.Range("Z10").Value = 9
$ws.Range("AG10").Value = 11

# Row 11 (Suwon FC vs Jeju SK) odds updates
$ws.Range("N11").Value = 2.03
$ws.Range("O11").Value = 1.78
$ws.Range("AJ11").Value = 29

# Row 13 (Trelleborg vs Helsingborg) odds updates
$ws.Range("N13").Value = 1.88
$ws.Range("O13").Value = 1.93
$ws.Range("R13").Value = 1.67

# Row 14 (Grasshoppers vs Aarau) odds updates
$ws.Range("G14").Value = 1.91
$ws.Range("H14").Value = 3.6
$ws.Range("I14").Value = 3.6
$ws.Range("J14").Value = 1.04
$ws.Range("K14").Value = 13
$ws.Range("N14").Value = 1.7
$ws.Range("O14").Value = 2.1
$ws.Range("R14").Value = 1.62
$ws.Range("S14").Value = 2.2
$ws.Range("T14").Value = 9
$ws.Range("U14").Value = 10
$ws.Range("W14").Value = 17
$ws.Range("X14").Value = 15
$ws.Range("AA14").Value = 7
$ws.Range("AE14").Value = 13
$ws.Range("AF14").Value = 21
$ws.Range("AG14").Value = 13
$ws.Range("AH14").Value = 41
$ws.Range("AI14").Value = 26
$ws.Range("AJ14").Value = 29
